{"js": "// Old division problem -> new division problem, as text replacements.\nconst pairs = [\n    [\"445\u00f73=\", \"190\u00f73=\"],\n    [\"874\u00f76=\", \"556\u00f78=\"],\n    [\"153\u00f76=\", \"682\u00f74=\"],\n    [\"236\u00f78=\", \"623\u00f78=\"],\n    [\"984\u00f72=\", \"546\u00f73=\"],\n    [\"291\u00f74=\", \"563\u00f79=\"],\n    [\"523\u00f78=\", \"480\u00f76=\"],\n    [\"552\u00f75=\", \"107\u00f77=\"],\n    [\"236\u00f79=\", \"115\u00f79=\"],\n    [\"907\u00f76=\", \"307\u00f79=\"],\n    [\"657\u00f72=\", \"990\u00f73=\"],\n    [\"941\u00f78=\", \"744\u00f76=\"],\n    [\"501\u00f79=\", \"518\u00f74=\"],\n    [\"531\u00f79=\", \"126\u00f78=\"],\n    [\"800\u00f77=\", \"144\u00f73=\"],\n    [\"346\u00f78=\", \"430\u00f72=\"],\n    [\"372\u00f77=\", \"979\u00f79=\"],\n    [\"587\u00f78=\", \"682\u00f72=\"],\n    [\"236\u00f73=\", \"327\u00f76=\"],\n    [\"973\u00f78=\", \"769\u00f75=\"],\n    [\"270\u00f78=\", \"522\u00f78=\"],\n    [\"566\u00f73=\", \"715\u00f77=\"],\n    [\"169\u00f73=\", \"777\u00f72=\"],\n    [\"797\u00f78=\", \"119\u00f77=\"],\n    [\"287\u00f75=\", \"885\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Old division problem -> new division problem, as text replacements.\n$pairs = @(\n    @(\"445\u00f73=\", \"190\u00f73=\"),\n    @(\"874\u00f76=\", \"556\u00f78=\"),\n    @(\"153\u00f76=\", \"682\u00f74=\"),\n    @(\"236\u00f78=\", \"623\u00f78=\"),\n    @(\"984\u00f72=\", \"546\u00f73=\"),\n    @(\"291\u00f74=\", \"563\u00f79=\"),\n    @(\"523\u00f78=\", \"480\u00f76=\"),\n    @(\"552\u00f75=\", \"107\u00f77=\"),\n    @(\"236\u00f79=\", \"115\u00f79=\"),\n    @(\"907\u00f76=\", \"307\u00f79=\"),\n    @(\"657\u00f72=\", \"990\u00f73=\"),\n    @(\"941\u00f78=\", \"744\u00f76=\"),\n    @(\"501\u00f79=\", \"518\u00f74=\"),\n    @(\"531\u00f79=\", \"126\u00f78=\"),\n    @(\"800\u00f77=\", \"144\u00f73=\"),\n    @(\"346\u00f78=\", \"430\u00f72=\"),\n    @(\"372\u00f77=\", \"979\u00f79=\"),\n    @(\"587\u00f78=\", \"682\u00f72=\"),\n    @(\"236\u00f73=\", \"327\u00f76=\"),\n    @(\"973\u00f78=\", \"769\u00f75=\"),\n    @(\"270\u00f78=\", \"522\u00f78=\"),\n    @(\"566\u00f73=\", \"715\u00f77=\"),\n    @(\"169\u00f73=\", \"777\u00f72=\"),\n    @(\"797\u00f78=\", \"119\u00f77=\"),\n    @(\"287\u00f75=\", \"885\u00f79=\"),\n)\n\nforeach ($pair in $pairs) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
